# The source CSV that feeds this sheet had a blank leading line and a few
# trailing blank lines. After cleaning the input, the header now sits on
# row 1 (was row 2) and the data block ends at row 44 (was row 45), with
# the old trailing blank rows 46-48 gone entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray blank first row - this shifts every row (header +
# data) up by one.
$ws.Rows.Item(1).Delete() | Out-Null

# Remove the now-trailing blank rows that used to be rows 46-48 (now
# rows 45-47 after the shift above) so the sheet's used range ends at
# row 44.
$ws.Range("A45:F47").EntireRow.Delete() | Out-Null

# Reflect where the author's cursor ended up after the cleanup.
$ws.Range("C46").Select() | Out-Null
